$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal TEXT value without letting Excel
# auto-convert numeric-looking strings (e.g. "1.003", "311.80") into
# real numbers. A leading apostrophe forces text entry (Excel's usual
# "quote prefix" convention) while the apostrophe itself is not stored
# as part of the value. The cells original Style is captured/restored
# so this is purely a value edit with no formatting side effect.
function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.Value = "'" + $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "28.008.61"
Set-TextValue $ws.Range("D3") "1.859.65"
Set-TextValue $ws.Range("E3") "  -0.64%  "
Set-TextValue $ws.Range("D4") "1.003"
Set-TextValue $ws.Range("E4") "  +0.27%  "
Set-TextValue $ws.Range("D5") "311.80"
Set-TextValue $ws.Range("E5") "  -0.24%  "
Set-TextValue $ws.Range("D6") "1.002"
Set-TextValue $ws.Range("E6") "  +0.21%  "
Set-TextValue $ws.Range("D7") "0.5088"
Set-TextValue $ws.Range("E7") "  +1.46%  "
Set-TextValue $ws.Range("D8") "0.3817"
Set-TextValue $ws.Range("E8") "  -0.26%  "
Set-TextValue $ws.Range("D9") "0.08308"
Set-TextValue $ws.Range("E9") "  -6.91%  "
Set-TextValue $ws.Range("D10") "1.113"
Set-TextValue $ws.Range("E10") "  -0.33%  "
Set-TextValue $ws.Range("E11") "  +0.21%  "
Set-TextValue $ws.Range("E12") "  -2.78%  "
Set-TextValue $ws.Range("D13") "20.53"
Set-TextValue $ws.Range("E13") "  -0.59%  "
Set-TextValue $ws.Range("D14") "1.852.34"
Set-TextValue $ws.Range("E14") "  -0.52%  "
Set-TextValue $ws.Range("D15") "7.206"
Set-TextValue $ws.Range("E15") "  -0.40%  "
Set-TextValue $ws.Range("E16") "  +0.29%  "
Set-TextValue $ws.Range("D17") "0.00001097"
Set-TextValue $ws.Range("E17") "  -0.09%  "
Set-TextValue $ws.Range("D18") "90.67"
Set-TextValue $ws.Range("E18") "  -0.47%  "
Set-TextValue $ws.Range("D19") "0.06625"
Set-TextValue $ws.Range("E19") "  -0.44%  "
Set-TextValue $ws.Range("D20") "17.66"
Set-TextValue $ws.Range("E20") "  -2.37%  "
Set-TextValue $ws.Range("E21") "  +0.15%  "
Set-TextValue $ws.Range("E22") "  -1.54%  "
Set-TextValue $ws.Range("D23") "28.018.73"
Set-TextValue $ws.Range("E23") "  +0.15%  "
Set-TextValue $ws.Range("D24") "11.05"
Set-TextValue $ws.Range("E24") "  -3.95%  "
Set-TextValue $ws.Range("D25") "2.239"
Set-TextValue $ws.Range("E25") "  -1.44%  "
Set-TextValue $ws.Range("D26") "2.544"
Set-TextValue $ws.Range("E26") "  +1.95%  "
Set-TextValue $ws.Range("D27") "2.076.08"
Set-TextValue $ws.Range("E27") "  -0.32%  "
Set-TextValue $ws.Range("D28") "157.99"
Set-TextValue $ws.Range("E28") "  -0.19%  "
Set-TextValue $ws.Range("D29") "20.47"
Set-TextValue $ws.Range("E29") "  -0.89%  "
Set-TextValue $ws.Range("D30") "124.52"
Set-TextValue $ws.Range("E30") "  -1.24%  "
Set-TextValue $ws.Range("D31") "0.1053"
Set-TextValue $ws.Range("E31") "  -0.69%  "
Set-TextValue $ws.Range("D32") "1.037"
Set-TextValue $ws.Range("E32") "  -1.69%  "
Set-TextValue $ws.Range("D33") "5.663"
Set-TextValue $ws.Range("E33") "  +1.16%  "
Set-TextValue $ws.Range("E34") "  -0.25%  "
Set-TextValue $ws.Range("D35") "9.406"
Set-TextValue $ws.Range("E35") "  -1.01%  "
Set-TextValue $ws.Range("B36") "Hedera"
Set-TextValue $ws.Range("C36") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D36") "0.06518"
Set-TextValue $ws.Range("E36") "  -0.47%  "
Set-TextValue $ws.Range("B37") "VeChain"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D37") "0.02410"
Set-TextValue $ws.Range("E37") "  +0.54%  "
Set-TextValue $ws.Range("D38") "0.2169"
Set-TextValue $ws.Range("E38") "  -0.55%  "
Set-TextValue $ws.Range("D39") "1.206"
Set-TextValue $ws.Range("E39") "  +0.23%  "
Set-TextValue $ws.Range("D40") "0.6441"
Set-TextValue $ws.Range("E40") "  +1.16%  "
Set-TextValue $ws.Range("D41") "1.219"
Set-TextValue $ws.Range("E41") "  -5.19%  "
Set-TextValue $ws.Range("D42") "4.897"
Set-TextValue $ws.Range("E42") "  -0.43%  "
Set-TextValue $ws.Range("D43") "11.19"
Set-TextValue $ws.Range("E43") "  -2.64%  "
Set-TextValue $ws.Range("D44") "0.6083"
Set-TextValue $ws.Range("E44") "  +1.31%  "
Set-TextValue $ws.Range("D45") "13.13"
Set-TextValue $ws.Range("E45") "  -0.23%  "
Set-TextValue $ws.Range("E46") "  +0.04%  "
Set-TextValue $ws.Range("D47") "3.657"
Set-TextValue $ws.Range("E47") "  -0.28%  "
Set-TextValue $ws.Range("D48") "2.016"
Set-TextValue $ws.Range("E48") "  +1.15%  "
Set-TextValue $ws.Range("D49") "1.208"
Set-TextValue $ws.Range("E49") "  -1.57%  "
Set-TextValue $ws.Range("D50") "120.25"
Set-TextValue $ws.Range("E50") "  -0.49%  "
Set-TextValue $ws.Range("D51") "78.68"
Set-TextValue $ws.Range("E51") "  -0.59%  "
